# Added one more Registration test (RegistrationWithDifferentPasswords)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RegistrationUser")

# Append the new test row (row 5): TestName, Email, FullName, Password, ConfirmPassword
$ws.Range("A5").Value = "RegistrationWithDifferentPasswords"
$ws.Range("B5").Value = "email@abv.bg"
$ws.Range("C5").Value = "softuni"
$ws.Range("D5").Value = "qa"
$ws.Range("E5").Value = "qa123"

# Make RegistrationUser the active/selected sheet with D9 selected
$ws.Select() | Out-Null
$ws.Range("D9").Select() | Out-Null
